# Edit: process.pptx slide with MGGPU diagram.
#  - Remove the solid (tx1) outline from four existing rounded-rectangle
#    shapes, replacing it with "no line".
#  - Nudge shape "Rounded Rectangle 158" (id 159) slightly.
#  - Add two new rounded-rectangle shapes (duplicated from the existing
#    style so fill / line / style / text-body all match) positioned over
#    the diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- Turn off the outline ("no line") on the four existing boxes ---
$shape154 = Get-ShapeById $s 154   # Rectangle: Rounded Corners 153
$shape154.Line.Visible = 0

$shape155 = Get-ShapeById $s 155   # Rectangle: Rounded Corners 154
$shape155.Line.Visible = 0

$shape159 = Get-ShapeById $s 159   # Rectangle: Rounded Corners 158
$shape159.Left = 851.42578125
$shape159.Top = 146.6641082763672
$shape159.Line.Visible = 0

$shape160 = Get-ShapeById $s 160   # Rectangle: Rounded Corners 159
$shape160.Line.Visible = 0

# --- Add two new rounded rectangles, matching the existing shape style ---
$newShape1 = $shape160.Duplicate()
$newShape1.Left = 770.0263061523438
$newShape1.Top = 213.80355834960938
$newShape1.Width = 146.00851440429688
$newShape1.Height = 79.06480407714844
$newShape1.Line.Visible = 0

$newShape2 = $shape160.Duplicate()
$newShape2.Left = 593.0526123046875
$newShape2.Top = 188.1372528076172
$newShape2.Width = 95.15220642089844
$newShape2.Height = 92.2293701171875
$newShape2.Line.Visible = 0
